$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "301.19"
Set-TextValue $ws "E2" "-3.07%"
Set-TextValue $ws "G2" "19"

# Row 3
Set-TextValue $ws "D3" "35.48"
Set-TextValue $ws "E3" "-0.28%"
Set-TextValue $ws "G3" "19"

# Row 4
Set-TextValue $ws "D4" "5.059"
Set-TextValue $ws "E4" "-1.15%"
Set-TextValue $ws "G4" "19"

# Row 5
Set-TextValue $ws "D5" "0.08009"
Set-TextValue $ws "E5" "-2.37%"
Set-TextValue $ws "G5" "19"

# Row 6
Set-TextValue $ws "E6" "-7.81%"
Set-TextValue $ws "G6" "19"

# Row 7
Set-TextValue $ws "D7" "7.759"
Set-TextValue $ws "E7" "-2.32%"
Set-TextValue $ws "G7" "19"

# Row 8
Set-TextValue $ws "D8" "0.9281"
Set-TextValue $ws "E8" "0.42%"
Set-TextValue $ws "G8" "19"

# Row 9
Set-TextValue $ws "D9" "0.1481"
Set-TextValue $ws "E9" "33.06%"
Set-TextValue $ws "G9" "19"

# Row 10
Set-TextValue $ws "D10" "0.1896"
Set-TextValue $ws "E10" "-1.06%"
Set-TextValue $ws "G10" "19"

# Row 11
Set-TextValue $ws "D11" "0.08981"
Set-TextValue $ws "E11" "-2.06%"
Set-TextValue $ws "G11" "19"

# Row 12
Set-TextValue $ws "D12" "0.03456"
Set-TextValue $ws "E12" "-5.52%"
Set-TextValue $ws "G12" "19"

# Row 13
Set-TextValue $ws "D13" "0.09868"
Set-TextValue $ws "E13" "-0.50%"
Set-TextValue $ws "G13" "19"

# Row 14
Set-TextValue $ws "D14" "0.001399"
Set-TextValue $ws "E14" "-3.07%"
Set-TextValue $ws "G14" "19"

# Row 15
Set-TextValue $ws "D15" "0.005784"
Set-TextValue $ws "E15" "-1.27%"
Set-TextValue $ws "G15" "19"

# Row 16
Set-TextValue $ws "E16" "1.59%"
Set-TextValue $ws "G16" "19"

# Row 17
Set-TextValue $ws "D17" "4.046"
Set-TextValue $ws "E17" "-1.99%"
Set-TextValue $ws "G17" "19"

# Row 18
Set-TextValue $ws "D18" "2.934"
Set-TextValue $ws "E18" "-1.85%"
Set-TextValue $ws "G18" "19"

# Row 19
Set-TextValue $ws "E19" "1.25%"
Set-TextValue $ws "G19" "19"

# Row 20
Set-TextValue $ws "D20" "0.1302"
Set-TextValue $ws "E20" "-0.57%"
Set-TextValue $ws "G20" "19"

# Row 21
Set-TextValue $ws "D21" "5.061"
Set-TextValue $ws "E21" "-0.73%"
Set-TextValue $ws "G21" "19"

# Row 22
Set-TextValue $ws "D22" "0.2395"
Set-TextValue $ws "E22" "8.60%"
Set-TextValue $ws "G22" "19"

# Row 23
Set-TextValue $ws "D23" "0.04487"
Set-TextValue $ws "E23" "-1.16%"
Set-TextValue $ws "G23" "19"

# Row 24
Set-TextValue $ws "D24" "0.001212"
Set-TextValue $ws "E24" "-1.16%"
Set-TextValue $ws "G24" "19"

# Row 25
Set-TextValue $ws "D25" "0.004763"
Set-TextValue $ws "E25" "-0.98%"
Set-TextValue $ws "G25" "19"

# Row 26
Set-TextValue $ws "D26" "0.0001228"
Set-TextValue $ws "E26" "-1.77%"
Set-TextValue $ws "G26" "19"

# Row 27
Set-TextValue $ws "D27" "0.0003019"
Set-TextValue $ws "E27" "-32.13%"
Set-TextValue $ws "G27" "19"

# Row 28
Set-TextValue $ws "G28" "19"

# Row 29
Set-TextValue $ws "G29" "19"

# Row 30
Set-TextValue $ws "G30" "19"

# Row 31
Set-TextValue $ws "G31" "19"

# Row 32
Set-TextValue $ws "G32" "19"

# Row 33
Set-TextValue $ws "G33" "19"

# Row 34
Set-TextValue $ws "G34" "19"

# Row 35
Set-TextValue $ws "G35" "19"

# Row 36
Set-TextValue $ws "G36" "19"

# Row 37
Set-TextValue $ws "G37" "19"

# Row 38
Set-TextValue $ws "G38" "19"

# Row 39
Set-TextValue $ws "D39" "0.01841"
Set-TextValue $ws "E39" "-7.00%"
Set-TextValue $ws "G39" "19"

# Row 40
Set-TextValue $ws "D40" "0.04765"
Set-TextValue $ws "E40" "-2.62%"
Set-TextValue $ws "G40" "19"

# Row 41
Set-TextValue $ws "E41" "9.96%"
Set-TextValue $ws "G41" "19"

# Row 42
Set-TextValue $ws "D42" "0.007313"
Set-TextValue $ws "E42" "-4.14%"
Set-TextValue $ws "G42" "19"

# Row 43
Set-TextValue $ws "E43" "-4.13%"
Set-TextValue $ws "G43" "19"

# Row 44
Set-TextValue $ws "D44" "0.002106"
Set-TextValue $ws "E44" "-2.41%"
Set-TextValue $ws "G44" "19"

# Row 45
Set-TextValue $ws "D45" "0.01082"
Set-TextValue $ws "E45" "-4.65%"
Set-TextValue $ws "G45" "19"

# Row 46
Set-TextValue $ws "D46" "0.00006229"
Set-TextValue $ws "E46" "-4.75%"
Set-TextValue $ws "G46" "19"

# Row 47
Set-TextValue $ws "E47" "-0.17%"
Set-TextValue $ws "G47" "19"

# Row 48
Set-TextValue $ws "D48" "64.66"
Set-TextValue $ws "E48" "-64.06%"
Set-TextValue $ws "G48" "19"

# Row 49
Set-TextValue $ws "G49" "19"

# Row 50
Set-TextValue $ws "D50" "0.00002097"
Set-TextValue $ws "E50" "-0.17%"
Set-TextValue $ws "G50" "19"

# Row 51
Set-TextValue $ws "D51" "0.0001997"
Set-TextValue $ws "E51" "-0.17%"
Set-TextValue $ws "G51" "19"
